$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.262.71"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +5.35%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.323.79"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.99%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "406.96"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.21%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "110.18"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.71%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.586"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +5.39%  "

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.05%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.635"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.51%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.57"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.84%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0985"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +5.86%  "

# Row 12
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.19%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.838.40"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.83%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.40"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.91%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "19.29"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.70%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.323.15"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.15%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.04"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.26%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "59.125.28"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +5.34%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.68"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -4.52%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.29"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.92%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0000108"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +4.04%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.83"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.91%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "306.82"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.96%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "75.14"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.22%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.21"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.70%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "28.56"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.49%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.45"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.23%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.79"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -4.70%  "

# Row 29
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.170"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.88%  "

# Row 30
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.25"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.17%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.04%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.112"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.39%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.34"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.44%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "40.39"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +9.70%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0522"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +6.49%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.13"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.06%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.85"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.81%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.25"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +5.00%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.05%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.45"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.68%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "137.98"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.22%  "

# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.91%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.88"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.19%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.90"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.39%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.61"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -5.09%  "

# Row 46
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.31"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +10.79%  "

# Row 47
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "TheGraph"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.275"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -3.18%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.16"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.42%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.171.25"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.97%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.44"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.49%  "

# Row 51
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.34"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +6.09%  "
